$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell value while forcing it to stay as plain text,
# so numeric-looking strings (e.g. "6.97", "0.670", "3.00") are not
# auto-converted into numbers by Excel (which would lose trailing
# zeros / exact formatting), and without leaving a permanent style
# change on the cell (NumberFormat is reset back to Normal afterward).
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.834.01"
Set-TextValue "E2" "  +3.10%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.946.69"
Set-TextValue "E3" "  +1.38%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "577.49"
Set-TextValue "E5" "  -0.49%  "

# Row 6 - Solana
Set-TextValue "D6" "150.30"
Set-TextValue "E6" "  +2.43%  "

# Row 7 - USDC
Set-TextValue "E7" "  +0.10%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "2.942.53"
Set-TextValue "E8" "  +1.30%  "

# Row 9 - XRP
Set-TextValue "E9" "  +0.10%  "

# Row 10 - Toncoin
Set-TextValue "D10" "6.97"

# Row 11 - Dogecoin
Set-TextValue "E11" "  -0.70%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.442"
Set-TextValue "E12" "  +1.92%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +0.68%  "

# Row 14 - Avalanche
Set-TextValue "D14" "33.97"
Set-TextValue "E14" "  +3.95%  "

# Row 15 - TRON
Set-TextValue "E15" "  +0.49%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue "D16" "3.433.30"
Set-TextValue "E16" "  +1.37%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "63.700.50"
Set-TextValue "E17" "  +2.95%  "

# Row 18 - Polkadot
Set-TextValue "D18" "6.83"
Set-TextValue "E18" "  +2.23%  "

# Row 19 - WrappedEther
Set-TextValue "D19" "2.943.30"
Set-TextValue "E19" "  +0.90%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "449.06"
Set-TextValue "E20" "  +2.82%  "

# Row 21 - Chainlink
Set-TextValue "D21" "13.52"
Set-TextValue "E21" "  +1.31%  "

# Row 22 - Polygon
Set-TextValue "D22" "0.670"
Set-TextValue "E22" "  +1.47%  "

# Row 23 - Uniswap
Set-TextValue "D23" "7.05"
Set-TextValue "E23" "  +1.34%  "

# Row 24 - Litecoin
Set-TextValue "D24" "79.94"
Set-TextValue "E24" "  -0.07%  "

# Row 25 - RenderToken
Set-TextValue "D25" "10.76"
Set-TextValue "E25" "  +5.11%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "12.12"
Set-TextValue "E26" "  +0.82%  "

# Row 27 - Fetch.AI
Set-TextValue "E27" "  +5.41%  "

# Row 29 - NEARProtocol
Set-TextValue "D29" "7.48"
Set-TextValue "E29" "  +4.99%  "

# Row 30 - PEPE
Set-TextValue "D30" "0.0000107"
Set-TextValue "E30" "  -3.77%  "

# Row 31 - PancakeSwap
Set-TextValue "E31" "  -0.23%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "2.12"
Set-TextValue "E32" "  +0.24%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.108"
Set-TextValue "E33" "  +0.47%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "26.28"
Set-TextValue "E34" "  +1.84%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "E35" "  -0.10%  "

# Row 36 - Mantle
Set-TextValue "D36" "0.967"
Set-TextValue "E36" "  +0.11%  "

# Row 37 - Stacks
Set-TextValue "E37" "  +7.05%  "

# Row 38 - Filecoin
Set-TextValue "E38" "  +0.94%  "

# Row 39 - dogwifhat
Set-TextValue "D39" "3.00"
Set-TextValue "E39" "  -2.60%  "

# Row 40 - OKB
Set-TextValue "D40" "49.04"
Set-TextValue "E40" "  -0.19%  "

# Row 41 - Arweave
Set-TextValue "D41" "43.83"
Set-TextValue "E41" "  +14.43%  "

# Row 42 - Kaspa
Set-TextValue "E42" "  +1.49%  "

# Row 43 - Cosmos
Set-TextValue "D43" "8.22"
Set-TextValue "E43" "  -1.31%  "

# Row 44 - TheGraph
Set-TextValue "E44" "  +4.95%  "

# Rows 45-47 rotate: VeChain/Maker/Bittensor -> Bittensor/VeChain/Maker
# Row 45 becomes Bittensor
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D45" "375.68"
Set-TextValue "E45" "  +9.20%  "

# Row 46 becomes VeChain
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0349"
Set-TextValue "E46" "  +3.93%  "

# Row 47 becomes Maker
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D47" "2.737.31"
Set-TextValue "E47" "  +1.85%  "

# Row 48 - Monero
Set-TextValue "D48" "133.97"
Set-TextValue "E48" "  -0.57%  "

# Row 50 - FLOKI
Set-TextValue "D50" "0.000218"
Set-TextValue "E50" "  +7.21%  "

# Row 51 - Stellar
Set-TextValue "D51" "0.105"
Set-TextValue "E51" "  +1.40%  "
